# Update column G ("K") values for rows 3-14 per the regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = 2
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 2
    10 = 4
    11 = 5
    12 = 3
    13 = 1
    14 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
